$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename HO_VectorBox -> HO_CAN across the Location column (D3:D34)
$rng = $ws.Range("D3:D34")
foreach ($cell in $rng.Cells) {
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "HO_VectorBox*") {
        $cell.Value = $val -replace "HO_VectorBox", "HO_CAN"
    }
}

# Add a new initial fault entry occupying rows 18 and 19
$ws.Range("C18").Value = "Output"
$ws.Range("D18").Value = "HO_CAN/Faults"
$ws.Range("C19").Value = "Output"
$ws.Range("D19").Value = "HO_CAN/Faults"

# Update the selected cell to C20
$ws.Range("C20").Select()
